# Auto-generated edit script: update cached market-price derived values
# per Sheets/Excalibur_Profits.xlsx scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 156.42857
$ws.Range("I2").Value = 165.83333
$ws.Range("K2").Value = 165.83333
$ws.Range("M2").Value = -52.83332999999999
$ws.Range("H19").Value = 699.46155
$ws.Range("I19").Value = 900.5
$ws.Range("K19").Value = 900.5
$ws.Range("M19").Value = -725.5
$ws.Range("H80").Value = 8929231
$ws.Range("I80").Value = 11364071
$ws.Range("J80").Value = 1483.3334
$ws.Range("K80").Value = 34092213
$ws.Range("L80").Value = 4450.0002
$ws.Range("M80").Value = -34091215
$ws.Range("N80").Value = -6446.0002
$ws.Range("H83").Value = 8929231
$ws.Range("I83").Value = 11364071
$ws.Range("J83").Value = 1483.3334
$ws.Range("K83").Value = 102276639
$ws.Range("L83").Value = 13350.0006
$ws.Range("M83").Value = -102271647
$ws.Range("N83").Value = -23334.0006
$ws.Range("H98").Value = 1325.7241
$ws.Range("I98").Value = 1341.9231
$ws.Range("K98").Value = 1341.9231
$ws.Range("M98").Value = 156.0769
$ws.Range("H100").Value = 3627.6206
$ws.Range("I100").Value = 1293.5714
$ws.Range("J100").Value = 9754.5
$ws.Range("K100").Value = 1293.5714
$ws.Range("L100").Value = 9754.5
$ws.Range("M100").Value = -752.5714
$ws.Range("N100").Value = -10836.5
$ws.Range("H116").Value = 29704
$ws.Range("I116").Value = 36574.668
$ws.Range("K116").Value = 36574.668
$ws.Range("M116").Value = -33132.668
$ws.Range("H122").Value = 1325.7241
$ws.Range("I122").Value = 1341.9231
$ws.Range("K122").Value = 4025.7693
$ws.Range("M122").Value = -1575.7693
$ws.Range("H132").Value = 31384.205
$ws.Range("I132").Value = 34431.543
$ws.Range("K132").Value = 103294.629
$ws.Range("M132").Value = -100764.629

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5888062
$ws.Range("I32").Value = 6333130
$ws.Range("J32").Value = 27999
$ws.Range("K32").Value = 6333130
$ws.Range("L32").Value = 27999
$ws.Range("M32").Value = -6332843
$ws.Range("N32").Value = -28573
$ws.Range("H37").Value = 45055.5
$ws.Range("I37").Value = 9000
$ws.Range("J37").Value = 67999.91
$ws.Range("K37").Value = 9000
$ws.Range("L37").Value = 67999.91
$ws.Range("M37").Value = -8727
$ws.Range("N37").Value = -68545.91
$ws.Range("H45").Value = 8012
$ws.Range("I45").Value = 5527.625
$ws.Range("K45").Value = 5527.625
$ws.Range("M45").Value = -5150.625
$ws.Range("H61").Value = 13143.5
$ws.Range("I61").Value = 11510.125
$ws.Range("K61").Value = 11510.125
$ws.Range("M61").Value = -11298.125
$ws.Range("H63").Value = 9874.875
$ws.Range("I63").Value = 1999.5
$ws.Range("K63").Value = 1999.5
$ws.Range("M63").Value = -1313.5
$ws.Range("H66").Value = 9874.875
$ws.Range("I66").Value = 1999.5
$ws.Range("K66").Value = 9997.5
$ws.Range("M66").Value = -6565.5
$ws.Range("H74").Value = 2577.8362
$ws.Range("J74").Value = 3515.9
$ws.Range("L74").Value = 3515.9
$ws.Range("N74").Value = -5263.9
$ws.Range("H77").Value = 2577.8362
$ws.Range("J77").Value = 3515.9
$ws.Range("L77").Value = 17579.5
$ws.Range("N77").Value = -26315.5
$ws.Range("H122").Value = 3837.5386
$ws.Range("I122").Value = 2498
$ws.Range("K122").Value = 7494
$ws.Range("M122").Value = -5044
$ws.Range("H136").Value = 13143.5
$ws.Range("I136").Value = 11510.125
$ws.Range("K136").Value = 34530.375
$ws.Range("M136").Value = -31980.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1353.3019
$ws.Range("I107").Value = 1378.6459
$ws.Range("K107").Value = 1378.6459
$ws.Range("M107").Value = 541.3541

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8564.223
$ws.Range("J31").Value = 11280.24
$ws.Range("L31").Value = 11280.24
$ws.Range("N31").Value = -11870.24
$ws.Range("H34").Value = 8564.223
$ws.Range("J34").Value = 11280.24
$ws.Range("L34").Value = 11280.24
$ws.Range("N34").Value = -11684.24
$ws.Range("H138").Value = 38802.6
$ws.Range("I138").Value = 39042.332
$ws.Range("J138").Value = 38699.855
$ws.Range("K138").Value = 39042.332
$ws.Range("L138").Value = 38699.855
$ws.Range("M138").Value = -33902.332
$ws.Range("N138").Value = -48979.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2974.75
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2974.75
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 4998.4
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 4998
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 14994
$ws.Range("M80").Value = -14064
$ws.Range("N80").Value = -16866
$ws.Range("H83").Value = 4998.4
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 4998
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 44982
$ws.Range("M83").Value = -40320
$ws.Range("N83").Value = -54342
$ws.Range("H107").Value = 617.875
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 591.8570999999999
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 1775.5713
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -5615.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1828.3939
$ws.Range("I97").Value = 303.6111
$ws.Range("K97").Value = 303.6111
$ws.Range("M97").Value = 192.3889
$ws.Range("H113").Value = 1835.8605
$ws.Range("I113").Value = 878.62067
$ws.Range("K113").Value = 878.62067
$ws.Range("M113").Value = 1291.37933
$ws.Range("H126").Value = 5832
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 9997
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 29991
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -34931

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12918.375
$ws.Range("I16").Value = 14707
$ws.Range("K16").Value = 14707
$ws.Range("M16").Value = -14537
$ws.Range("H22").Value = 47567.91
$ws.Range("I22").Value = 100875.4
$ws.Range("K22").Value = 100875.4
$ws.Range("M22").Value = -100580.4
$ws.Range("H27").Value = 47567.91
$ws.Range("I27").Value = 100875.4
$ws.Range("K27").Value = 100875.4
$ws.Range("M27").Value = -100768.4
$ws.Range("H29").Value = 2000000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
